$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 26.144619
$ws.Range("H2").Value = 78.433857
$ws.Range("I2").Value = 0.5211737020083955
$ws.Range("J2").Value = 0.5211737020083955
$ws.Range("M2").Value = 2.166102
$ws.Range("N2").Value = 6.498306
$ws.Range("O2").Value = 0.04231752823769151
$ws.Range("P2").Value = 0.04231752823769151
$ws.Range("Q2").Value = 56.63191150513801
$ws.Range("R2").Value = 509.6872035462421
$ws.Range("S2").Value = 0.02205478285148249
$ws.Range("T2").Value = 0.0220547828514825
$ws.Range("G3").Value = 26.144619
$ws.Range("H3").Value = 78.433857
$ws.Range("I3").Value = 0.5211737020083955
$ws.Range("J3").Value = 0.5211737020083955
$ws.Range("O3").Value = 0.5792223765593866
$ws.Range("P3").Value = 0.5792223765593866
$ws.Range("Q3").Value = 775.150906424877
$ws.Range("R3").Value = 6976.358157823893
$ws.Range("S3").Value = 0.3018754702775564
$ws.Range("T3").Value = 0.3018754702775564
$ws.Range("G4").Value = 26.144619
$ws.Range("H4").Value = 78.433857
$ws.Range("I4").Value = 0.5211737020083955
$ws.Range("J4").Value = 0.5211737020083955
$ws.Range("M4").Value = 19.37218933333333
$ws.Range("N4").Value = 58.116568
$ws.Range("O4").Value = 0.3784600952029219
$ws.Range("P4").Value = 0.378460095202922
$ws.Range("Q4").Value = 506.4785093158641
$ws.Range("R4").Value = 4558.306583842776
$ws.Range("S4").Value = 0.1972434488793566
$ws.Range("T4").Value = 0.1972434488793567
$ws.Range("I5").Value = 0.3571392594830743
$ws.Range("J5").Value = 0.3571392594830742
$ws.Range("M5").Value = 2.166102
$ws.Range("N5").Value = 6.498306
$ws.Range("O5").Value = 0.04231752823769151
$ws.Range("P5").Value = 0.04231752823769151
$ws.Range("Q5").Value = 38.80755851670001
$ws.Range("R5").Value = 349.2680266503
$ws.Range("S5").Value = 0.01511325069796323
$ws.Range("T5").Value = 0.01511325069796323
$ws.Range("I6").Value = 0.3571392594830743
$ws.Range("J6").Value = 0.3571392594830742
$ws.Range("O6").Value = 0.5792223765593866
$ws.Range("P6").Value = 0.5792223765593866
$ws.Range("Q6").Value = 531.1795657405501
$ws.Range("S6").Value = 0.2068630506404457
$ws.Range("T6").Value = 0.2068630506404457
$ws.Range("I7").Value = 0.3571392594830743
$ws.Range("J7").Value = 0.3571392594830742
$ws.Range("M7").Value = 19.37218933333333
$ws.Range("N7").Value = 58.116568
$ws.Range("O7").Value = 0.3784600952029219
$ws.Range("P7").Value = 0.378460095202922
$ws.Range("Q7").Value = 347.0692382676001
$ws.Range("R7").Value = 3123.6231444084
$ws.Range("S7").Value = 0.1351629581446653
$ws.Range("T7").Value = 0.1351629581446653
$ws.Range("G8").Value = 6.104416333333333
$ws.Range("H8").Value = 18.313249
$ws.Range("I8").Value = 0.1216870385085301
$ws.Range("J8").Value = 0.1216870385085301
$ws.Range("M8").Value = 2.166102
$ws.Range("N8").Value = 6.498306
$ws.Range("O8").Value = 0.04231752823769151
$ws.Range("P8").Value = 0.04231752823769151
$ws.Range("Q8").Value = 13.222788428466
$ws.Range("R8").Value = 119.005095856194
$ws.Range("S8").Value = 0.005149494688245778
$ws.Range("T8").Value = 0.005149494688245779
$ws.Range("G9").Value = 6.104416333333333
$ws.Range("H9").Value = 18.313249
$ws.Range("I9").Value = 0.1216870385085301
$ws.Range("J9").Value = 0.1216870385085301
$ws.Range("O9").Value = 0.5792223765593866
$ws.Range("P9").Value = 0.5792223765593866
$ws.Range("Q9").Value = 180.987294325389
$ws.Range("R9").Value = 1628.885648928501
$ws.Range("S9").Value = 0.07048385564138443
$ws.Range("T9").Value = 0.07048385564138443
$ws.Range("G10").Value = 6.104416333333333
$ws.Range("H10").Value = 18.313249
$ws.Range("I10").Value = 0.1216870385085301
$ws.Range("J10").Value = 0.1216870385085301
$ws.Range("M10").Value = 19.37218933333333
$ws.Range("N10").Value = 58.116568
$ws.Range("O10").Value = 0.3784600952029219
$ws.Range("P10").Value = 0.378460095202922
$ws.Range("Q10").Value = 118.2559089788258
$ws.Range("R10").Value = 1064.303180809432
$ws.Range("S10").Value = 0.04605368817889994
$ws.Range("T10").Value = 0.04605368817889995
